$d = $word.ActiveDocument

# Locate the unique anchor text "new_address[0] }}." which ends the sentence
# "... to {{ new_address[0] }}."
$find = $d.Content.Find
$find.ClearFormatting()
$found = $find.Execute(
    "new_address[0] }}.",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "",
    0
)

if (-not $found) {
    throw "Could not locate 'new_address[0] }}.' anchor text"
}

$full = $find.Parent.Duplicate

# The trailing run's text is " }}." (space, brace, brace, period) -- 4 chars.
# Remove it entirely, then re-insert just the "." so it becomes its own run
# (kept distinct from the preceding "[0]" run) before appending the new
# "{{ new_address[0].city }}, {{ new_address[0].state }}" pieces.
$tailStart = $full.End - 4
$tail = $d.Range($tailStart, $full.End)
$tail.Delete()

$insertPoint = $d.Range($tailStart, $tailStart)
$insertPoint.InsertAfter(".")
$insertPoint.Collapse(0)
$r0 = $d.Range($tailStart, $insertPoint.Start)
$r0.Font.Name = "Times New Roman"
$r0.Font.Size = 12

# Now append four more runs:
#   "city"
#   " }}"
#   ", {{ new_address[0].state }}"
#   "."
$insertPoint.InsertAfter("city")
$insertPoint.Collapse(0)
$r1 = $d.Range($insertPoint.Start - 4, $insertPoint.Start)
$r1.Font.Name = "Times New Roman"
$r1.Font.Size = 12

$insertPoint.InsertAfter(" }}")
$insertPoint.Collapse(0)
$r2 = $d.Range($insertPoint.Start - 3, $insertPoint.Start)
$r2.Font.Name = "Times New Roman"
$r2.Font.Size = 12

$insertPoint.InsertAfter(", {{ new_address[0].state }}")
$insertPoint.Collapse(0)
$r3 = $d.Range($insertPoint.Start - 29, $insertPoint.Start)
$r3.Font.Name = "Times New Roman"
$r3.Font.Size = 12

$insertPoint.InsertAfter(".")
$insertPoint.Collapse(0)
$r4 = $d.Range($insertPoint.Start - 1, $insertPoint.Start)
$r4.Font.Name = "Times New Roman"
$r4.Font.Size = 12

Write-Output "done"
